$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 108
$ws.Range("F6").Value = 1765
$ws.Range("F7").Value = 938
$ws.Range("F8").Value = 584
$ws.Range("F9").Value = 2724
$ws.Range("F10").Value = 755
$ws.Range("F11").Value = 580
$ws.Range("F13").Value = 57
$ws.Range("F15").Value = 363
$ws.Range("F16").Value = 370
$ws.Range("F18").Value = 2138
$ws.Range("F20").Value = 721
$ws.Range("F22").Value = 2638
$ws.Range("F28").Value = 509
$ws.Range("F32").Value = 556
$ws.Range("F33").Value = 218
$ws.Range("F35").Value = 358
$ws.Range("F37").Value = 186

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 4222
$ws.Range("F10").Value = 322
$ws.Range("F25").Value = 18
$ws.Range("F26").Value = 284
$ws.Range("F33").Value = 25
$ws.Range("F35").Value = 489
$ws.Range("F36").Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1439
$ws.Range("F6").Value = 555
$ws.Range("F7").Value = 172
$ws.Range("F8").Value = 221

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1439
$ws.Range("F4").Value = 555
$ws.Range("F5").Value = 108
$ws.Range("F9").Value = 1765
$ws.Range("F10").Value = 172
$ws.Range("F11").Value = 938
$ws.Range("F12").Value = 584
$ws.Range("F13").Value = 2725
$ws.Range("F14").Value = 755
$ws.Range("F15").Value = 580
$ws.Range("F17").Value = 57
$ws.Range("F19").Value = 363
$ws.Range("F21").Value = 370
$ws.Range("F25").Value = 2138
$ws.Range("F27").Value = 721
$ws.Range("F30").Value = 2638
$ws.Range("F37").Value = 221
$ws.Range("F39").Value = 509
$ws.Range("F40").Value = 509
$ws.Range("F43").Value = 556
$ws.Range("F44").Value = 284
$ws.Range("F45").Value = 218
$ws.Range("F47").Value = 358
$ws.Range("F49").Value = 186
$ws.Range("F50").Value = 489
